# Regenerate column G ("K") values in the save_data sheet.
# The original column used a "Strike#" based computation; it is being
# regenerated here (std/mean recalculated elsewhere) and the resulting
# values (s_vals) are written back into column G for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G (7)
$kValues = @{
    2 = 0
    3 = 1
    5 = 2
    6 = 2
    7 = 0
    8 = 2
    9 = 2
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 2
    16 = 3
    17 = 2
    18 = 1
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 4
    30 = 4
    31 = 0
    32 = 1
    33 = 2
    34 = 0
    35 = 0
    36 = 2
    37 = 1
    38 = 1
    39 = 1
    40 = 3
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 2
    46 = 1
    47 = 3
    48 = 2
    49 = 2
    50 = 2
    52 = 3
    54 = 1
    55 = 1
    56 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
